$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.147.00"
$ws.Range("D3").Value = "2.805.15"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "360.71"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.38"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.30"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.56"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.66"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "3.243.50"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "2.785.86"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.952"
$ws.Range("E17").Value = "  +7.03%  "
$ws.Range("D18").Value = "52.125.53"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.14"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "0.0₃0988"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.70"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.52"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  +5.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.52"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.59"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.96"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.38"
$ws.Range("E45").Value = "  -6.45%  "
$ws.Range("D46").Value = "2.084.85"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.73"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.934"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.95"
